$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to add to rows 2-5 (columns B..H are plain text, even numeric-looking
# values such as "12" or "1"; column A holds real numbers styled like the
# bold/bordered header row).
$data = @(
    @{ A = 0; B = "Oxak";      C = "кг"; D = "12"; E = "доллары"; F = "3";  G = "36"; H = "effrfwwpkp" },
    @{ A = 1; B = "Sjsnisb";   C = "кг"; D = "5";  E = "суммы";   F = "12"; G = "60"; H = "effrfwwpkp" },
    @{ A = 2; B = "111111111"; C = "kg"; D = "1";  E = "суммы";   F = "1";  G = "1";  H = "effrfwwpkp" },
    @{ A = 3; B = "22222222";  C = "кг"; D = "1";  E = "суммы";   F = "1";  G = "1";  H = "effrfwwpkp" }
)

# Force text storage for B:H up front so the numeric-looking strings ("12",
# "1", "36", ...) are kept as text instead of being auto-coerced to numbers.
$textRange = $ws.Range("B2:H5")
$textRange.NumberFormat = "@"

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.A

    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = $item.E
    $ws.Cells.Item($row, 6).Value = $item.F
    $ws.Cells.Item($row, 7).Value = $item.G
    $ws.Cells.Item($row, 8).Value = $item.H

    $row++
}

# Drop the temporary text-number-format so the cells fall back to the
# workbook's default (unstyled) look, matching the target cells which carry
# no style index.
$textRange.Style = "Normal"

# Give column A (rows 2-5) the same bold/border/center-top look as the header
# row by copying the header's cell format onto each of them.
$ws.Range("B1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
